# Update the "想去人数" (F column) figures in both the "展览" and "全部类型"
# worksheets to reflect the newly generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 8363
    3  = 7862
    5  = 192
    8  = 130
    10 = 170
    14 = 1819
    19 = 128
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
